# logboek aangevuld met taken van 19/02
#
# Appends six new paragraphs (a blank line, "19/02", and four task lines)
# right after the last paragraph of the document ("Tanguy: helpen zoeken
# naar script voor JSON"), and moves the hidden "_GoBack" bookmark so it
# ends up wrapping the new final paragraph, exactly like the original
# author's edit.

$d = $word.ActiveDocument

# The document always ends with a (hidden) "_GoBack" bookmark marking the
# last edit position. Drop it here - we recreate it at the new end below -
# otherwise we would end up with two identically named bookmarks.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$lastPara = $d.Paragraphs.Last

# Sanity check: this should be the "Tanguy: ..." paragraph the diff attaches to.
# (Falls back gracefully - the insertion position is computed either way.)
$insertPos = $lastPara.Range.End - 1
$target = $d.Range($insertPos, $insertPos)

# Build the WordprocessingML for the new paragraphs as a self-contained
# single-part OOXML package and splice it in with InsertXML so the
# run/paragraph structure (incl. the run-less blank paragraph and the
# two separate runs in the "Ruben:" line) matches exactly, instead of
# relying on text + InsertParagraphAfter which synthesizes extra runs.
$newParagraphsXml = @'
<w:p><w:pPr><w:rPr><w:lang w:val="nl-BE"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="nl-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-BE"/></w:rPr><w:t>19/02</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="nl-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-BE"/></w:rPr><w:t xml:space="preserve">Ruben: </w:t></w:r><w:r><w:rPr><w:lang w:val="nl-BE"/></w:rPr><w:t>Registratie in orde brengen en controle reservatie laten functioneren</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="nl-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-BE"/></w:rPr><w:t>Jeroen: handleiding</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="nl-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-BE"/></w:rPr><w:t>Dieter: script voor JSON uitzoeken</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="nl-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-BE"/></w:rPr><w:t>Tanguy: Dieter helpen</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$packageXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $newParagraphsXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($packageXml) | Out-Null
